$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing "Baseline Members" (previously "Foreigners") row data ---
$ws.Range("B12").Value = "Baseline Members"
$ws.Range("C12").Value = 134
$ws.Range("D12").Value = 60.91

$ws.Range("B13").Value = "Baseline Members"
$ws.Range("C13").Value = 86
$ws.Range("D13").Value = 39.09

# --- Add new "New Members" rows ---
$ws.Range("A14").Value = "Able to meet essential needs"
$ws.Range("B14").Value = "New Members"
$ws.Range("C14").Value = 36
$ws.Range("D14").Value = 28.35

$ws.Range("A15").Value = "Unable to meet essential needs"
$ws.Range("B15").Value = "New Members"
$ws.Range("C15").Value = 91
$ws.Range("D15").Value = 71.65

# --- Column widths (best effort match to bestFit widths of 29.71 / 17.71 / 4) ---
$ws.Columns.Item(1).ColumnWidth = 28.833333
$ws.Columns.Item(2).ColumnWidth = 16.833333
$ws.Columns.Item(3).ColumnWidth = 3.166667
